# chore: update Sheets via scheduled runner
# Refresh market-board derived figures (currentAveragePrice*, LevePrice*,
# LeveProfit* columns H:N) across the per-job sheets with the latest
# pulled prices.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3040.5
$ws.Range("J69").Value = 3010
$ws.Range("L69").Value = 9030
$ws.Range("N69").Value = -10778

$ws.Range("H72").Value = 3040.5
$ws.Range("J72").Value = 3010
$ws.Range("L72").Value = 27090
$ws.Range("N72").Value = -35826

$ws.Range("H132").Value = 3923041.8
$ws.Range("I132").Value = 4652113
$ws.Range("J132").Value = 4283.25
$ws.Range("K132").Value = 13956339
$ws.Range("L132").Value = 12849.75
$ws.Range("M132").Value = -13953809
$ws.Range("N132").Value = -17909.75

$ws.Range("H137").Value = 2183.0637
$ws.Range("I137").Value = 2150.147
$ws.Range("J137").Value = 2269.1538
$ws.Range("K137").Value = 6450.441
$ws.Range("L137").Value = 6807.4614
$ws.Range("M137").Value = -3900.441
$ws.Range("N137").Value = -11907.4614

$ws.Range("H141").Value = 954837.4
$ws.Range("I141").Value = 2197.1428
$ws.Range("K141").Value = 6591.428400000001
$ws.Range("M141").Value = -1411.428400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 793.9167
$ws.Range("I74").Value = 685.3929000000001
$ws.Range("J74").Value = 1173.75
$ws.Range("K74").Value = 685.3929000000001
$ws.Range("L74").Value = 1173.75
$ws.Range("M74").Value = 188.6070999999999
$ws.Range("N74").Value = -2921.75

$ws.Range("H77").Value = 793.9167
$ws.Range("I77").Value = 685.3929000000001
$ws.Range("J77").Value = 1173.75
$ws.Range("K77").Value = 3426.9645
$ws.Range("L77").Value = 5868.75
$ws.Range("M77").Value = 941.0355
$ws.Range("N77").Value = -14604.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3568.5715
$ws.Range("I20").Value = 1636
$ws.Range("J20").Value = 8400
$ws.Range("K20").Value = 1636
$ws.Range("L20").Value = 8400
$ws.Range("M20").Value = -1389
$ws.Range("N20").Value = -8894

$ws.Range("H94").Value = 2454.5
$ws.Range("I94").Value = 2009
$ws.Range("K94").Value = 2009
$ws.Range("M94").Value = -1558

$ws.Range("H134").Value = 2822.8333
$ws.Range("I134").Value = 1962.4
$ws.Range("K134").Value = 5887.200000000001
$ws.Range("M134").Value = -3352.200000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3537
$ws.Range("I16").Value = 2074
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 2074
$ws.Range("L16").Value = 5000
$ws.Range("M16").Value = -1787
$ws.Range("N16").Value = -5574

$ws.Range("H113").Value = 3537
$ws.Range("I113").Value = 2074
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 2074
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = 96
$ws.Range("N113").Value = -9340

$ws.Range("H132").Value = 4641.7646
$ws.Range("I132").Value = 3264
$ws.Range("J132").Value = 5866.4443
$ws.Range("K132").Value = 9792
$ws.Range("L132").Value = 17599.3329
$ws.Range("M132").Value = -7262
$ws.Range("N132").Value = -22659.3329

$ws.Range("H134").Value = 1655.0571
$ws.Range("I134").Value = 1239.8846
$ws.Range("J134").Value = 2854.4443
$ws.Range("K134").Value = 3719.6538
$ws.Range("L134").Value = 8563.332900000001
$ws.Range("M134").Value = -1184.6538
$ws.Range("N134").Value = -13633.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2300.6736
$ws.Range("I68").Value = 701.3158
$ws.Range("J68").Value = 3313.6
$ws.Range("K68").Value = 2103.9474
$ws.Range("L68").Value = 9940.799999999999
$ws.Range("M68").Value = -1292.9474
$ws.Range("N68").Value = -11562.8

$ws.Range("H71").Value = 2300.6736
$ws.Range("I71").Value = 701.3158
$ws.Range("J71").Value = 3313.6
$ws.Range("K71").Value = 6311.8422
$ws.Range("L71").Value = 29822.4
$ws.Range("M71").Value = -2255.8422
$ws.Range("N71").Value = -37934.39999999999

$ws.Range("H107").Value = 910.05084
$ws.Range("I107").Value = 543
$ws.Range("J107").Value = 1241.5807
$ws.Range("K107").Value = 1629
$ws.Range("L107").Value = 3724.7421
$ws.Range("M107").Value = 291
$ws.Range("N107").Value = -7564.742099999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 38626.547
$ws.Range("J42").Value = 38626.547
$ws.Range("L42").Value = 38626.547
$ws.Range("N42").Value = -39596.547

$ws.Range("H113").Value = 10000
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws.Range("H115").Value = 38626.547
$ws.Range("J115").Value = 38626.547
$ws.Range("L115").Value = 38626.547
$ws.Range("N115").Value = -40976.547

$ws.Range("H132").Value = 3528.6
$ws.Range("I132").Value = 2145.6667
$ws.Range("K132").Value = 6437.000100000001
$ws.Range("M132").Value = -3907.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()

$ws.Range("H122").Value = 3145.1614
$ws.Range("I122").Value = 2620
$ws.Range("J122").Value = 5333.3335
$ws.Range("K122").Value = 7860
$ws.Range("L122").Value = 16000.0005
$ws.Range("M122").Value = -5410
$ws.Range("N122").Value = -20900.0005

$ws.Range("H132").Value = 3401.9092
$ws.Range("I132").Value = 2262.6667
$ws.Range("J132").Value = 4769
$ws.Range("K132").Value = 6788.000100000001
$ws.Range("L132").Value = 14307
$ws.Range("M132").Value = -4258.000100000001
$ws.Range("N132").Value = -19367

$ws.Range("H136").Value = 1852.1666
$ws.Range("I136").Value = 1615.52
$ws.Range("K136").Value = 4846.559999999999
$ws.Range("M136").Value = -2296.559999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1185.3334
$ws.Range("I81").Value = 655
$ws.Range("J81").Value = 1450.5
$ws.Range("K81").Value = 1310
$ws.Range("L81").Value = 2901
$ws.Range("M81").Value = -249
$ws.Range("N81").Value = -5023

$ws.Range("H84").Value = 1185.3334
$ws.Range("I84").Value = 655
$ws.Range("J84").Value = 1450.5
$ws.Range("K84").Value = 6550
$ws.Range("L84").Value = 14505
$ws.Range("M84").Value = -1246
$ws.Range("N84").Value = -25113

$ws.Range("H107").Value = 2253.818
$ws.Range("I107").Value = 534.5714
$ws.Range("J107").Value = 5262.5
$ws.Range("K107").Value = 1603.7142
$ws.Range("L107").Value = 15787.5
$ws.Range("M107").Value = 316.2857999999999
$ws.Range("N107").Value = -19627.5

$ws.Range("H132").Value = 15632.632
$ws.Range("I132").Value = 1901.48
$ws.Range("J132").Value = 42038.69
$ws.Range("K132").Value = 5704.440000000001
$ws.Range("L132").Value = 126116.07
$ws.Range("M132").Value = -3174.440000000001
$ws.Range("N132").Value = -131176.07

$ws.Range("H136").Value = 2310.087
$ws.Range("I136").Value = 1183.2307
$ws.Range("J136").Value = 3775
$ws.Range("K136").Value = 3549.6921
$ws.Range("L136").Value = 11325
$ws.Range("M136").Value = -999.6921000000002
$ws.Range("N136").Value = -16425
